# Daily attendance processing - 2025-11-12 13:36:51
# Swap the first two comma-separated entries in the "Recorded By" column (G)
# for every data row, leaving any additional entries in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 157
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -ne "") {
        $parts = $text.Split(",")

        if ($parts.Length -ge 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()

            $rest = ""
            for ($i = 2; $i -lt $parts.Length; $i++) {
                $rest = $rest + ", " + $parts[$i].Trim()
            }

            $newVal = $second + ", " + $first + $rest

            if ($newVal -ne $text) {
                $cell.Value = $newVal
            }
        }
    }
}
